$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-8: farsnocc -> Log(TFVN), with model column uppercased
$ws.Range("A2:A8").Value = "Log(TFVN)"
$ws.Range("B2").Value = "OLS"
$ws.Range("B3").Value = "SAR"
$ws.Range("B4").Value = "SEM"
$ws.Range("B5").Value = "SARAR"
$ws.Range("B6").Value = "SDM"
$ws.Range("B7").Value = "SDEM"
$ws.Range("B8").Value = "SLX"

# Rows 9-15: farsocc -> Log(TFVO), with model column uppercased
$ws.Range("A9:A15").Value = "Log(TFVO)"
$ws.Range("B9").Value = "OLS"
$ws.Range("B10").Value = "SAR"
$ws.Range("B11").Value = "SEM"
$ws.Range("B12").Value = "SARAR"
$ws.Range("B13").Value = "SDM"
$ws.Range("B14").Value = "SDEM"
$ws.Range("B15").Value = "SLX"
